$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Timestamps for rows 2..75 (time_taken column, F)
$timeTaken = @(
    "2021-10-05 13:39:48.668749",
    "2021-10-05 13:39:48.668761",
    "2021-10-05 13:39:48.668765",
    "2021-10-05 13:39:48.668768",
    "2021-10-05 13:39:48.668771",
    "2021-10-05 13:39:48.668774",
    "2021-10-05 13:39:48.668776",
    "2021-10-05 13:39:48.668779",
    "2021-10-05 13:39:48.668782",
    "2021-10-05 13:39:48.668785",
    "2021-10-05 13:39:48.668788",
    "2021-10-05 13:39:48.668791",
    "2021-10-05 13:39:48.668793",
    "2021-10-05 13:39:48.668796",
    "2021-10-05 13:39:48.668799",
    "2021-10-05 13:39:48.668802",
    "2021-10-05 13:39:48.668805",
    "2021-10-05 13:39:48.668808",
    "2021-10-05 13:39:48.668810",
    "2021-10-05 13:39:48.668813",
    "2021-10-05 13:39:48.668816",
    "2021-10-05 13:39:48.668819",
    "2021-10-05 13:39:48.668821",
    "2021-10-05 13:39:48.668824",
    "2021-10-05 13:39:48.668830",
    "2021-10-05 13:39:48.668834",
    "2021-10-05 13:39:48.668839",
    "2021-10-05 13:39:48.668844",
    "2021-10-05 13:39:48.668848",
    "2021-10-05 13:39:48.668853",
    "2021-10-05 13:39:48.668857",
    "2021-10-05 13:39:48.668862",
    "2021-10-05 13:39:48.668893",
    "2021-10-05 13:39:48.668907",
    "2021-10-05 13:39:48.668936",
    "2021-10-05 13:39:48.668943",
    "2021-10-05 13:39:48.668946",
    "2021-10-05 13:39:48.668949",
    "2021-10-05 13:39:48.668952",
    "2021-10-05 13:39:48.668955",
    "2021-10-05 13:39:48.668959",
    "2021-10-05 13:39:48.668962",
    "2021-10-05 13:39:48.668965",
    "2021-10-05 13:39:48.668968",
    "2021-10-05 13:39:48.668970",
    "2021-10-05 13:39:48.668973",
    "2021-10-05 13:39:48.668976",
    "2021-10-05 13:39:48.668979",
    "2021-10-05 13:39:48.668982",
    "2021-10-05 13:39:48.668985",
    "2021-10-05 13:39:48.668988",
    "2021-10-05 13:39:48.668991",
    "2021-10-05 13:39:48.668994",
    "2021-10-05 13:39:48.668997",
    "2021-10-05 13:39:48.669000",
    "2021-10-05 13:39:48.669003",
    "2021-10-05 13:39:48.669006",
    "2021-10-05 13:39:48.669009",
    "2021-10-05 13:39:48.669012",
    "2021-10-05 13:39:48.669015",
    "2021-10-05 13:39:48.669017",
    "2021-10-05 13:39:48.669020",
    "2021-10-05 13:39:48.669023",
    "2021-10-05 13:39:48.669026",
    "2021-10-05 13:39:48.669031",
    "2021-10-05 13:39:48.669034",
    "2021-10-05 13:39:48.669036",
    "2021-10-05 13:39:48.669039",
    "2021-10-05 13:39:48.669042",
    "2021-10-05 13:39:48.669044",
    "2021-10-05 13:39:48.669047",
    "2021-10-05 13:39:48.669049",
    "2021-10-05 13:39:48.669052",
    "2021-10-05 13:39:48.669055"
)

# Copy formatting from the existing header cell (E1) onto the new header cell (F1)
# so that F1 gets the same bold/centered/bordered style used by the other headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F1").Value = "time_taken"

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
